$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix F646 value
$ws.Cells.Item(646, 6).Value = 5566.35789109

# Append new rows 647-657
$ws.Cells.Item(647, 1).Value = 45116.41666666666
$ws.Cells.Item(647, 2).Value = 30299.25
$ws.Cells.Item(647, 3).Value = 30453.27
$ws.Cells.Item(647, 4).Value = 30080.24
$ws.Cells.Item(647, 5).Value = 30174.62
$ws.Cells.Item(647, 6).Value = 5874.45936717
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(647, 1).PasteSpecial(-4122)

$ws.Cells.Item(648, 1).Value = 45117.41666666666
$ws.Cells.Item(648, 2).Value = 30175.34
$ws.Cells.Item(648, 3).Value = 31042.51
$ws.Cells.Item(648, 4).Value = 29965.03
$ws.Cells.Item(648, 5).Value = 30423.95
$ws.Cells.Item(648, 6).Value = 18369.45647798
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(648, 1).PasteSpecial(-4122)

$ws.Cells.Item(649, 1).Value = 45118.41666666666
$ws.Cells.Item(649, 2).Value = 30422.95
$ws.Cells.Item(649, 3).Value = 30809.56
$ws.Cells.Item(649, 4).Value = 30320.36
$ws.Cells.Item(649, 5).Value = 30631.36
$ws.Cells.Item(649, 6).Value = 14390.16504579
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(649, 1).PasteSpecial(-4122)

$ws.Cells.Item(650, 1).Value = 45119.41666666666
$ws.Cells.Item(650, 2).Value = 30633.89
$ws.Cells.Item(650, 3).Value = 30982
$ws.Cells.Item(650, 4).Value = 30227.25
$ws.Cells.Item(650, 5).Value = 30396.78
$ws.Cells.Item(650, 6).Value = 20184.77143358
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(650, 1).PasteSpecial(-4122)

$ws.Cells.Item(651, 1).Value = 45120.41666666666
$ws.Cells.Item(651, 2).Value = 30395.64
$ws.Cells.Item(651, 3).Value = 31829
$ws.Cells.Item(651, 4).Value = 30258.46
$ws.Cells.Item(651, 5).Value = 31482.21
$ws.Cells.Item(651, 6).Value = 36831.45497786
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(651, 1).PasteSpecial(-4122)

$ws.Cells.Item(652, 1).Value = 45121.41666666666
$ws.Cells.Item(652, 2).Value = 31483.23
$ws.Cells.Item(652, 3).Value = 31644.47
$ws.Cells.Item(652, 4).Value = 29940.08
$ws.Cells.Item(652, 5).Value = 30333.65
$ws.Cells.Item(652, 6).Value = 28504.11311169
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(652, 1).PasteSpecial(-4122)

$ws.Cells.Item(653, 1).Value = 45122.41666666666
$ws.Cells.Item(653, 2).Value = 30332.66
$ws.Cells.Item(653, 3).Value = 30403.97
$ws.Cells.Item(653, 4).Value = 30267.04
$ws.Cells.Item(653, 5).Value = 30299
$ws.Cells.Item(653, 6).Value = 4039.37147264
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(653, 1).PasteSpecial(-4122)

$ws.Cells.Item(654, 1).Value = 45123.41666666666
$ws.Cells.Item(654, 2).Value = 30300.6
$ws.Cells.Item(654, 3).Value = 30457.63
$ws.Cells.Item(654, 4).Value = 30078.23
$ws.Cells.Item(654, 5).Value = 30250.49
$ws.Cells.Item(654, 6).Value = 6357.2037676
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(654, 1).PasteSpecial(-4122)

$ws.Cells.Item(655, 1).Value = 45124.41666666666
$ws.Cells.Item(655, 2).Value = 30248.97
$ws.Cells.Item(655, 3).Value = 30342.59
$ws.Cells.Item(655, 4).Value = 29678.15
$ws.Cells.Item(655, 5).Value = 30154.32
$ws.Cells.Item(655, 6).Value = 16010.77083874
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(655, 1).PasteSpecial(-4122)

$ws.Cells.Item(656, 1).Value = 45125.41666666666
$ws.Cells.Item(656, 2).Value = 30152.07
$ws.Cells.Item(656, 3).Value = 30243.7
$ws.Cells.Item(656, 4).Value = 29522.25
$ws.Cells.Item(656, 5).Value = 29868.81
$ws.Cells.Item(656, 6).Value = 16104.96081001
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(656, 1).PasteSpecial(-4122)

$ws.Cells.Item(657, 1).Value = 45126.41666666666
$ws.Cells.Item(657, 2).Value = 29863.81
$ws.Cells.Item(657, 3).Value = 30201.29
$ws.Cells.Item(657, 4).Value = 29770.34
$ws.Cells.Item(657, 5).Value = 29921.83
$ws.Cells.Item(657, 6).Value = 12551.08599458
$ws.Cells.Item(646, 1).Copy()
$ws.Cells.Item(657, 1).PasteSpecial(-4122)

